$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()
$ws.Rows.Item(60).Insert()

$ws.Cells.Item(60, 1).Value = 5
$ws.Cells.Item(60, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(60, 3).Value = "Maule"
$ws.Cells.Item(60, 4).Value = 44895
$ws.Cells.Item(60, 5).Value = 7
$ws.Cells.Item(60, 6).Value = 100112022
$ws.Cells.Item(60, 7).Value = "Arveja Verde"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 500
$ws.Cells.Item(60, 11).Value = 19000
$ws.Cells.Item(60, 12).Value = 19000
$ws.Cells.Item(60, 13).Value = 19000
$ws.Cells.Item(60, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(60, 15).Value = "Región del Maule"
$ws.Cells.Item(60, 16).Value = 760
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"

$ws.Cells.Item(61, 1).Value = 5
$ws.Cells.Item(61, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(61, 3).Value = "Maule"
$ws.Cells.Item(61, 4).Value = 44895
$ws.Cells.Item(61, 5).Value = 7
$ws.Cells.Item(61, 6).Value = 100112022
$ws.Cells.Item(61, 7).Value = "Arveja Verde"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 400
$ws.Cells.Item(61, 11).Value = 19000
$ws.Cells.Item(61, 12).Value = 20000
$ws.Cells.Item(61, 13).Value = 19500
$ws.Cells.Item(61, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(61, 15).Value = "Región del Maule"
$ws.Cells.Item(61, 16).Value = 780
$ws.Cells.Item(61, 17).Value = 25
$ws.Cells.Item(61, 18).Value = "Hortaliza"

Write-Output "done"
